$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be written as text (matching the original
    # inline-string / shared-string storage) rather than letting Excel
    # auto-convert a numeric-looking string (e.g. "1.79e-02") into a number.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2 (Suffrage)
$ws.Range("C2").Value = 0.6752876796897646
$ws.Range("D2").Value = 1.009167506784137
$ws.Range("E2").Value = 0.02204064037765832
Set-TextValue $ws.Range("G2") "1.79e-02"

# Row 3 (GatesS)
$ws.Range("C3").Value = 0.5036067805454411
$ws.Range("D3").Value = 1.134293978280183
$ws.Range("E3").Value = 0.007913767507918164
$ws.Range("F3").Value = "norm_coldread_coverage_line_%"
Set-TextValue $ws.Range("G3") "9.00e-02"

# Row 4 (GatesT)
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0.4865644270916379
$ws.Range("D4").Value = 1.046600125451598
$ws.Range("E4").Value = -0.01688470457019641
$ws.Range("F4").Value = "norm_coldread_saccade_regression_rate_%"
Set-TextValue $ws.Range("G4") "3.39e-01"
